$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the duplicated client record in row 6 ---
# Row 6 had accidentally been filled with the same id_razonsocial /
# nombre_empresa as row 2. Replace with the correct, distinct values so the
# "processed / not processed" comparison logic (and the duplicate-value
# highlighting added below) works correctly.
$ws.Range("B6").Value = "30-70892538-7"
$ws.Range("C6").Value = "CONDOMINIO SIL FBSF SA"

# --- Add conditional formatting to flag duplicate values ---
# id_razonsocial (column B) and nombre_empresa (column C) each get a
# "Highlight Duplicate Values" rule across the full column, so repeated
# client records (like the one just fixed) are easy to spot going forward.
$rangeB = $ws.Range("B1:B1048576")
$fcB = $rangeB.FormatConditions.AddUniqueValues()
$fcB.DupeUnique = 1

$rangeC = $ws.Range("C1:C1048576")
$fcC = $rangeC.FormatConditions.AddUniqueValues()
$fcC.DupeUnique = 1

$fcB.Priority = 2
$fcC.Priority = 1

# Column C uses the standard "Light Red Fill with Dark Red Text" style.
$fcC.Font.Color = 393372
$fcC.Interior.Color = 13551615

# (A rule was also tried out on column D with the green style and then
# removed again while choosing the final look.)
$rangeD = $ws.Range("D1:D1048576")
$fcD = $rangeD.FormatConditions.AddUniqueValues()
$fcD.DupeUnique = 1
$fcD.Font.Color = 24832
$fcD.Interior.Color = 13561798
$rangeD.FormatConditions.Delete()

# Column B uses the "Green Fill with Dark Green Text" style.
$fcB.Font.Color = 24832
$fcB.Interior.Color = 13561798

# --- Selection / view state ---
# Leave the active cell on B6 and scroll back so row 1 is visible.
$ws.Range("B6").Select()
